# Adding .docx file generation feature
# -------------------------------------------------------------------
# Reproduces the two user-visible/content-level edits of the target
# workbook revision (the rest of the upstream diff is Excel-version
# metadata noise - new mc:/xr:/x14ac: namespaces, revisionPtr GUIDs,
# calcId/feature bumps, cosmetic theme-name strings, etc. - that gets
# rewritten automatically whenever a file is re-saved by a newer Excel
# build and isn't something a user action via the object model drives):
#
#   1. Column widths A:D were fitted to their content (AutoFit/Best Fit).
#   2. A page footer watermark ("RESTRICTED", orange Calibri 22pt) was
#      added to the right section of the sheet's footer.
# -------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column widths (A:D), matching the fitted widths recorded in the
#    saved workbook (A=21.71, B=4.57, C=73, D=7.86 characters).
$ws.Columns.Item(1).ColumnWidth = 21.7109375
$ws.Columns.Item(2).ColumnWidth = 4.5703125
$ws.Columns.Item(3).ColumnWidth = 73
$ws.Columns.Item(4).ColumnWidth = 7.85546875

# 2) Footer watermark: right-aligned "RESTRICTED" in orange (FF8939),
#    Calibri, 22pt, preceded by a carriage-return/format-reset code.
$ws.PageSetup.RightFooter = "`r&1#&`"Calibri`"&22&KFF8939 RESTRICTED"
